$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")

# Update column widths (B, C, E changed; target OOXML widths 56.7109375 / 51.7109375 / 56.7109375).
# The ColumnWidth COM setter quantizes to a 1/6-character pixel grid, so these inputs are chosen
# to land on the closest achievable grid value to the target width.
$ws.Columns.Item(2).ColumnWidth = 55.85
$ws.Columns.Item(3).ColumnWidth = 50.85
$ws.Columns.Item(5).ColumnWidth = 55.85

# Update timetable cell contents (B2:F13)
$ws.Range("B2").Value = '{}'
$ws.Range("C2").Value = '{0: sala nr 1 | Jan Nowak | Język polski}'
$ws.Range("D2").Value = '{}'
$ws.Range("E2").Value = '{}'
$ws.Range("F2").Value = '{}'
$ws.Range("B3").Value = '{}'
$ws.Range("C3").Value = '{0: sala nr 6 | Jan Nowak | Język polski}'
$ws.Range("D3").Value = '{}'
$ws.Range("E3").Value = '{}'
$ws.Range("F3").Value = '{}'
$ws.Range("B4").Value = '{}'
$ws.Range("C4").Value = '{0: sala nr 11 | Paweł Lewandowski | Matematyka}'
$ws.Range("D4").Value = '{}'
$ws.Range("E4").Value = '{}'
$ws.Range("F4").Value = '{}'
$ws.Range("B5").Value = '{0: sala nr 4 | Paweł Lewandowski | Matematyka}'
$ws.Range("C5").Value = '{0: sala nr 5 | Karolina Kamińska | Chemia}'
$ws.Range("D5").Value = '{}'
$ws.Range("E5").Value = '{}'
$ws.Range("F5").Value = '{0: sala nr 3 | Paweł Lewandowski | Matematyka}'
$ws.Range("B6").Value = '{0: sala nr 10 | Natalia Szymańska | Geografia}'
$ws.Range("C6").Value = '{0: sala nr 4 | Mateusz Kowalski | Język niemiecki}'
$ws.Range("D6").Value = '{}'
$ws.Range("E6").Value = '{0: sala nr 7 | Dominik Kaczor | Informatyka}'
$ws.Range("F6").Value = '{}'
$ws.Range("B7").Value = '{0: sala nr 10 | Zofia Wiśniewska | Wychowanie fizyczne}'
$ws.Range("C7").Value = '{0: sala nr 8 | Katarzyna Mazur | Fizyka}'
$ws.Range("D7").Value = '{0: sala nr 4 | Dominik Kaczor | Informatyka}'
$ws.Range("E7").Value = '{0: sala nr 11 | Lena Kowalska | Język angielski}'
$ws.Range("F7").Value = '{0: sala nr 10 | Katarzyna Mazur | Fizyka}'
$ws.Range("B8").Value = '{0: sala nr 8 | Paweł Lewandowski | Matematyka}'
$ws.Range("C8").Value = '{0: sala nr 2 | Natalia Szymańska | Geografia}'
$ws.Range("D8").Value = '{0: sala nr 3 | Dominik Kaczor | Informatyka}'
$ws.Range("E8").Value = '{0: sala nr 10 | Zofia Wiśniewska | Wychowanie fizyczne}'
$ws.Range("F8").Value = '{0: sala nr 9 | Piotr Wójcik | Biologia}'
$ws.Range("B9").Value = '{}'
$ws.Range("C9").Value = '{}'
$ws.Range("D9").Value = '{0: sala nr 7 | Dominik Kaczor | Informatyka}'
$ws.Range("E9").Value = '{0: sala nr 3 | Piotr Wójcik | Biologia}'
$ws.Range("F9").Value = '{0: sala nr 1 | Dominik Kaczor | Informatyka}'
$ws.Range("B10").Value = '{}'
$ws.Range("C10").Value = '{}'
$ws.Range("D10").Value = '{0: sala nr 8 | Lena Kowalska | Język angielski}'
$ws.Range("E10").Value = '{}'
$ws.Range("F10").Value = '{0: sala nr 8 | Jan Nowak | Język polski}'
$ws.Range("B11").Value = '{}'
$ws.Range("C11").Value = '{}'
$ws.Range("D11").Value = '{0: sala nr 2 | Mateusz Kowalski | Język niemiecki}'
$ws.Range("E11").Value = '{}'
$ws.Range("F11").Value = '{0: sala nr 7 | Paweł Lewandowski | Matematyka}'
$ws.Range("B12").Value = '{}'
$ws.Range("C12").Value = '{}'
$ws.Range("D12").Value = '{0: sala nr 7 | Karolina Kamińska | Chemia}'
$ws.Range("E12").Value = '{}'
$ws.Range("F12").Value = '{0: sala nr 3 | Zofia Wiśniewska | Wychowanie fizyczne}'
$ws.Range("B13").Value = '{}'
$ws.Range("C13").Value = '{}'
$ws.Range("D13").Value = '{0: sala nr 11 | Katarzyna Mazur | Fizyka}'
$ws.Range("E13").Value = '{}'
$ws.Range("F13").Value = '{0: sala nr 11 | Lena Kowalska | Język angielski}'
